$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; everything shifts one column to the right
$ws.Columns("A:A").Insert()

# Fill in the new column B (previously column A) with the new "local code" values
$ws.Range("B2").Value = 301
$ws.Range("B3").Value = 302
$ws.Range("B4").Value = 303

# Update the descriptive text columns (previously B and D, now C and E)
$ws.Range("C2").Value = "ACADEMIA TESTE"
$ws.Range("E2").Value = "ACADEMIA TESTE"
$ws.Range("C3").Value = "ADM TESTE"
$ws.Range("E3").Value = "ADM TESTE"
$ws.Range("C4").Value = "ALMOXARIFADOS TESTE"
$ws.Range("E4").Value = "ALMOXARIFADOS TESTE"

# Fill in the new first column (EMPRESA)
$ws.Range("A1").Value = "EMPRESA"
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 1

# Give the underline formatting to D3 (RESPONSAVEL cell on the ADM TESTE row)
$ws.Range("D3").Font.Underline = $true

# Resize columns to match the new layout spacing
$ws.Columns("A:A").ColumnWidth = 14.6
$ws.Columns("C:C").ColumnWidth = 21.3

# Move the active selection to A4
$ws.Range("A4").Select()
